# Apply the commit: add "Sprint" / "Due Date" columns (M, N) to the
# "Result" sheet, update their row values/formatting, move the selection
# on the "URL" sheet from D11 to D7, and touch the workbook view/calc
# metadata.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Result" sheet — add columns M (Sprint) and N (Due Date)
# ---------------------------------------------------------------------
$result = $wb.Worksheets.Item("Result")

$result.Range("M1").Value = "Sprint"
$result.Range("N1").Value = "Due Date"

# Rows 2-12: M gets the same "Passed" styling used elsewhere (copy format
# from an existing "Passed" cell so the shared fill style is reused
# instead of minting a new one), N gets the "needs attention" yellow
# styling used by the other "Update/Add ..." cells in the sheet.
for ($row = 2; $row -le 12; $row++) {
    $passedSrc = $result.Range("B" + $row)
    $passedSrc.Copy()
    $mCell = $result.Range("M" + $row)
    $mCell.PasteSpecial(-4122)
    $mCell.Value = "Passed"

    $flagSrc = $result.Range("K" + $row)
    $flagSrc.Copy()
    $nCell = $result.Range("N" + $row)
    $nCell.PasteSpecial(-4122)
    $nCell.Value = "Add Due Date"
}

[void]$result.Range("A1").Select()

# ---------------------------------------------------------------------
# "URL" sheet — move the active selection from D11 to D7
# ---------------------------------------------------------------------
$url = $wb.Worksheets.Item("URL")
$url.Activate()
[void]$url.Range("D7").Select()
